# Updated cryptos list on Sat Aug 17 02:53:21 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table, and fixes the ShibaInu / WrappedEther rows which were in the
# wrong order (row 17 should be WrappedEther, row 18 should be ShibaInu).
#
# Values that look numeric (e.g. "521.65", "0.999") are written through a
# helper that forces the Text number format first and clears the format
# afterwards, so Excel keeps them as literal strings (matching how the
# source data is stored) instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "59.003.89"
Set-TextValue $ws.Range("E2") "  +1.93%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.589.28"
Set-TextValue $ws.Range("E3") "  +0.65%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "521.65"
Set-TextValue $ws.Range("E5") "  +1.02%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "139.84"
Set-TextValue $ws.Range("E6") "  -2.17%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  +0.05%  "

# Row 8 - XRP
Set-TextValue $ws.Range("E8") "  -0.20%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.600.57"
Set-TextValue $ws.Range("E9") "  +0.48%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("E10") "  -0.84%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("E11") "  +0.18%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("E12") "  +1.91%  "

# Row 13 - TRON
Set-TextValue $ws.Range("E13") "  +3.17%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.043.70"
Set-TextValue $ws.Range("E14") "  +0.53%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "59.083.23"
Set-TextValue $ws.Range("E15") "  +2.06%  "

# Row 16 - Avalanche
Set-TextValue $ws.Range("D16") "20.44"
Set-TextValue $ws.Range("E16") "  +0.81%  "

# Rows 17/18 - ShibaInu and WrappedEther were swapped; row 17 now holds
# WrappedEther's (refreshed) data, row 18 now holds ShibaInu's.
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "2.611.48"
Set-TextValue $ws.Range("E17") "  +3.33%  "

Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.0000133"
Set-TextValue $ws.Range("E18") "  -0.30%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "338.88"
Set-TextValue $ws.Range("E19") "  -0.07%  "

# Row 20 - Polkadot
Set-TextValue $ws.Range("E20") "  +0.29%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("E21") "  -0.90%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("E22") "  +3.10%  "

# Row 23 - Dai
Set-TextValue $ws.Range("D23") "0.999"

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "66.19"
Set-TextValue $ws.Range("E24") "  +1.11%  "

# Row 25 - Kaspa
Set-TextValue $ws.Range("E25") "  +1.02%  "

# Row 26 - Polygon
Set-TextValue $ws.Range("E26") "  +0.56%  "

# Row 27 - Binance-PegBSC-USD
Set-TextValue $ws.Range("E27") "  -0.36%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("E28") "  +0.70%  "

# Row 29 - USDe
Set-TextValue $ws.Range("D29") "0.999"

# Row 30 - PEPE
Set-TextValue $ws.Range("D30") "0.0₃0725"
Set-TextValue $ws.Range("E30") "  -3.26%  "

# Row 31 - Aptos
Set-TextValue $ws.Range("D31") "5.96"
Set-TextValue $ws.Range("E31") "  -5.23%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("E32") "  -0.08%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("E33") "  +0.28%  "

# Row 34 - Monero
Set-TextValue $ws.Range("D34") "148.88"
Set-TextValue $ws.Range("E34") "  -0.38%  "

# Row 35 - NEARProtocol
Set-TextValue $ws.Range("E35") "  -0.23%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("E36") "  -1.86%  "

# Row 37 - OKB
Set-TextValue $ws.Range("D37") "36.79"
Set-TextValue $ws.Range("E37") "  +1.79%  "

# Row 38 - Stacks
Set-TextValue $ws.Range("E38") "  +0.62%  "

# Row 39 - Fetch.AI
Set-TextValue $ws.Range("E39") "  -0.76%  "

# Row 40 - SuiNetwork
Set-TextValue $ws.Range("E40") "  -7.05%  "

# Row 41 - Filecoin
Set-TextValue $ws.Range("D41") "3.50"
Set-TextValue $ws.Range("E41") "  -0.60%  "

# Row 42 - FirstDigitalUSD
Set-TextValue $ws.Range("E42") "  +0.02%  "

# Row 43 - Bittensor
Set-TextValue $ws.Range("D43") "272.79"
Set-TextValue $ws.Range("E43") "  +0.66%  "

# Row 44 - WhiteBITCoin
Set-TextValue $ws.Range("E44") "  +1.05%  "

# Row 45 - Mantle
Set-TextValue $ws.Range("D45") "0.590"
Set-TextValue $ws.Range("E45") "  +0.50%  "

# Row 46 - Stellar
Set-TextValue $ws.Range("E46") "  -0.35%  "

# Row 47 - Hedera
Set-TextValue $ws.Range("E47") "  -0.69%  "

# Row 48 - EnergySwap
Set-TextValue $ws.Range("E48") "  -1.62%  "

# Row 49 - Maker
Set-TextValue $ws.Range("D49") "1.972.49"
Set-TextValue $ws.Range("E49") "  -0.17%  "

# Row 50 - RenderToken
Set-TextValue $ws.Range("D50") "4.65"
Set-TextValue $ws.Range("E50") "  +1.08%  "

# Row 51 - VeChain
Set-TextValue $ws.Range("E51") "  -0.26%  "
